$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# New row 152 : 2012-03-16, 粤APN645, 公司 -> 家, mileage 156714 -> 156829
# ---------------------------------------------------------------------------
$ws.Cells.Item(152, 1).Value = 40984
$ws.Cells.Item(152, 2).Value = "17:45-21:00"
$ws.Cells.Item(152, 3).Value = "粤APN645"
$ws.Cells.Item(152, 4).Value = "公司"
$ws.Cells.Item(152, 5).Value = "家"
$ws.Cells.Item(152, 6).Value = 156714
$ws.Cells.Item(152, 7).Value = 156829
$ws.Cells.Item(152, 8).Formula = "=G152-F152"

# match the date / time-range formatting already used by the rows above
$ws.Cells.Item(151, 2).Copy()
$ws.Cells.Item(152, 2).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# New row 153 : 2012-03-19, 粤APN645, 家 -> 公司, mileage 156848 -> 156955
# ---------------------------------------------------------------------------
$ws.Cells.Item(153, 1).Value = 40987
$ws.Cells.Item(153, 2).Value = "5:40-8:10"
$ws.Cells.Item(153, 3).Value = "粤APN645"
$ws.Cells.Item(153, 4).Value = "家"
$ws.Cells.Item(153, 5).Value = "公司"
$ws.Cells.Item(153, 6).Value = 156848
$ws.Cells.Item(153, 7).Value = 156955
$ws.Cells.Item(153, 8).Formula = "=G153-F153"

$ws.Cells.Item(151, 2).Copy()
$ws.Cells.Item(153, 2).PasteSpecial(-4122)

# note in J153 -- text begins with '+' so Excel stores it with a quote-prefix
$ws.Cells.Item(153, 10).Value = "'+油 ￥40"

# ---------------------------------------------------------------------------
# New note in J146 (existing row)
# ---------------------------------------------------------------------------
$ws.Cells.Item(146, 10).Value = "违章 -200"

# ---------------------------------------------------------------------------
# Move the view: active cell / selection on the frozen bottom-right pane
# ---------------------------------------------------------------------------
$ws.Range("G157").Select()
